$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1: A1/B1 get new FAQ/response values, C1/D1/E1 stay the same
$ws.Range("A1").Value = "faq1e"
$ws.Range("B1").Value = "resp1e"
$ws.Range("C1").Value = "admission"
$ws.Range("D1").Value = "rules"
$ws.Range("E1").Value = "payment"

# Update row 2: A2/B2 get new FAQ/response values, C2 stays the same
# (set B2 before A2 so the shared-strings table order matches the target)
$ws.Range("B2").Value = "res2e"
$ws.Range("A2").Value = "faq2e"
$ws.Range("C2").Value = "rules"

# Update the active selection to D2
$ws.Range("D2").Select()
